$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update cell text values (shared strings updated to 2025 FE data) ---
$ws.Range("B2").Value = "Jul 2024 - Jun 2025"
$ws.Range("E2").Value = "<p>`n  ONS published a `n  <a href=`"https://osr.statisticsauthority.gov.uk/correspondence/michael-keoghan-to-siobhan-tuohy-smith-request-to-suspend-aps-accreditation/`">response to OSR</a> about the current quality of Annual Population Survey (APS) (and Labour Force Survey) outputs. ONS asked OSR to temporarily suspend accreditation of all APS-based ONS outputs. There has since been a `n  <a href=`"https://osr.statisticsauthority.gov.uk/correspondence/ed-humpherson-to-michael-keoghan-suspension-of-the-accredited-official-statistics-status-for-the-estimates-ons-produces-from-the-annual-population-survey/`">response letter from OSR</a>. Overall, ONS’ view on the quality of the APS is that while it is robust for national and headline regional estimates, there are concerns with the quality of estimates for smaller segments of the population, such as local authority geographies. ONS will publish an explanatory note later this year providing guidance to users on the quality of current APS and will be used to inform further work ONS is undertaking to improve quality of the survey.`n</p>"
$ws.Range("B3").Value = "Jul 2024 - Jun 2025"
$ws.Range("E3").Value = "<p>`n  ONS published a `n  <a href=`"https://osr.statisticsauthority.gov.uk/correspondence/michael-keoghan-to-siobhan-tuohy-smith-request-to-suspend-aps-accreditation/`">response to OSR</a> about the current quality of Annual Population Survey (APS) (and Labour Force Survey) outputs. ONS asked OSR to temporarily suspend accreditation of all APS-based ONS outputs. There has since been a `n  <a href=`"https://osr.statisticsauthority.gov.uk/correspondence/ed-humpherson-to-michael-keoghan-suspension-of-the-accredited-official-statistics-status-for-the-estimates-ons-produces-from-the-annual-population-survey/`">response letter from OSR</a>. Overall, ONS’ view on the quality of the APS is that while it is robust for national and headline regional estimates, there are concerns with the quality of estimates for smaller segments of the population, such as local authority geographies. ONS will publish an explanatory note later this year providing guidance to users on the quality of current APS and will be used to inform further work ONS is undertaking to improve quality of the survey.`n</p>"
$ws.Range("B4").Value = "Jul 2024 - Jun 2025"
$ws.Range("E4").Value = "<p>`n  ONS published a `n  <a href=`"https://osr.statisticsauthority.gov.uk/correspondence/michael-keoghan-to-siobhan-tuohy-smith-request-to-suspend-aps-accreditation/`">response to OSR</a> about the current quality of Annual Population Survey (APS) (and Labour Force Survey) outputs. ONS asked OSR to temporarily suspend accreditation of all APS-based ONS outputs. There has since been a `n  <a href=`"https://osr.statisticsauthority.gov.uk/correspondence/ed-humpherson-to-michael-keoghan-suspension-of-the-accredited-official-statistics-status-for-the-estimates-ons-produces-from-the-annual-population-survey/`">response letter from OSR</a>. Overall, ONS’ view on the quality of the APS is that while it is robust for national and headline regional estimates, there are concerns with the quality of estimates for smaller segments of the population, such as local authority geographies. ONS will publish an explanatory note later this year providing guidance to users on the quality of current APS and will be used to inform further work ONS is undertaking to improve quality of the survey.`n</p>"
$ws.Range("B5").Value = "Jul 2024 - Jun 2025"
$ws.Range("E5").Value = "<p>`n  ONS published a `n  <a href=`"https://osr.statisticsauthority.gov.uk/correspondence/michael-keoghan-to-siobhan-tuohy-smith-request-to-suspend-aps-accreditation/`">response to OSR</a> about the current quality of Annual Population Survey (APS) (and Labour Force Survey) outputs. ONS asked OSR to temporarily suspend accreditation of all APS-based ONS outputs. There has since been a `n  <a href=`"https://osr.statisticsauthority.gov.uk/correspondence/ed-humpherson-to-michael-keoghan-suspension-of-the-accredited-official-statistics-status-for-the-estimates-ons-produces-from-the-annual-population-survey/`">response letter from OSR</a>. Overall, ONS’ view on the quality of the APS is that while it is robust for national and headline regional estimates, there are concerns with the quality of estimates for smaller segments of the population, such as local authority geographies. ONS will publish an explanatory note later this year providing guidance to users on the quality of current APS and will be used to inform further work ONS is undertaking to improve quality of the survey.`n</p>"
$ws.Range("B6").Value = "Jul 2024 - Jun 2025"
$ws.Range("E6").Value = "<p>`n  ONS published a `n  <a href=`"https://osr.statisticsauthority.gov.uk/correspondence/michael-keoghan-to-siobhan-tuohy-smith-request-to-suspend-aps-accreditation/`">response to OSR</a> about the current quality of Annual Population Survey (APS) (and Labour Force Survey) outputs. ONS asked OSR to temporarily suspend accreditation of all APS-based ONS outputs. There has since been a `n  <a href=`"https://osr.statisticsauthority.gov.uk/correspondence/ed-humpherson-to-michael-keoghan-suspension-of-the-accredited-official-statistics-status-for-the-estimates-ons-produces-from-the-annual-population-survey/`">response letter from OSR</a>. Overall, ONS’ view on the quality of the APS is that while it is robust for national and headline regional estimates, there are concerns with the quality of estimates for smaller segments of the population, such as local authority geographies. ONS will publish an explanatory note later this year providing guidance to users on the quality of current APS and will be used to inform further work ONS is undertaking to improve quality of the survey.`n</p>"
$ws.Range("B7").Value = "Jul 2024 - Jun 2025"
$ws.Range("E7").Value = "<p>`n  ONS published a `n  <a href=`"https://osr.statisticsauthority.gov.uk/correspondence/michael-keoghan-to-siobhan-tuohy-smith-request-to-suspend-aps-accreditation/`">response to OSR</a> about the current quality of Annual Population Survey (APS) (and Labour Force Survey) outputs. ONS asked OSR to temporarily suspend accreditation of all APS-based ONS outputs. There has since been a `n  <a href=`"https://osr.statisticsauthority.gov.uk/correspondence/ed-humpherson-to-michael-keoghan-suspension-of-the-accredited-official-statistics-status-for-the-estimates-ons-produces-from-the-annual-population-survey/`">response letter from OSR</a>. Overall, ONS’ view on the quality of the APS is that while it is robust for national and headline regional estimates, there are concerns with the quality of estimates for smaller segments of the population, such as local authority geographies. ONS will publish an explanatory note later this year providing guidance to users on the quality of current APS and will be used to inform further work ONS is undertaking to improve quality of the survey.`n</p>"
$ws.Range("B8").Value = "Jul 2024 - Jun 2025"
$ws.Range("E8").Value = "<p>`n  ONS published a `n  <a href=`"https://osr.statisticsauthority.gov.uk/correspondence/michael-keoghan-to-siobhan-tuohy-smith-request-to-suspend-aps-accreditation/`">response to OSR</a> about the current quality of Annual Population Survey (APS) (and Labour Force Survey) outputs. ONS asked OSR to temporarily suspend accreditation of all APS-based ONS outputs. There has since been a `n  <a href=`"https://osr.statisticsauthority.gov.uk/correspondence/ed-humpherson-to-michael-keoghan-suspension-of-the-accredited-official-statistics-status-for-the-estimates-ons-produces-from-the-annual-population-survey/`">response letter from OSR</a>. Overall, ONS’ view on the quality of the APS is that while it is robust for national and headline regional estimates, there are concerns with the quality of estimates for smaller segments of the population, such as local authority geographies. ONS will publish an explanatory note later this year providing guidance to users on the quality of current APS and will be used to inform further work ONS is undertaking to improve quality of the survey.`n</p>"
$ws.Range("B9").Value = "Jul 2024 - Jun 2025"
$ws.Range("E9").Value = "<p>`n  ONS published a `n  <a href=`"https://osr.statisticsauthority.gov.uk/correspondence/michael-keoghan-to-siobhan-tuohy-smith-request-to-suspend-aps-accreditation/`">response to OSR</a> about the current quality of Annual Population Survey (APS) (and Labour Force Survey) outputs. ONS asked OSR to temporarily suspend accreditation of all APS-based ONS outputs. There has since been a `n  <a href=`"https://osr.statisticsauthority.gov.uk/correspondence/ed-humpherson-to-michael-keoghan-suspension-of-the-accredited-official-statistics-status-for-the-estimates-ons-produces-from-the-annual-population-survey/`">response letter from OSR</a>. Overall, ONS’ view on the quality of the APS is that while it is robust for national and headline regional estimates, there are concerns with the quality of estimates for smaller segments of the population, such as local authority geographies. ONS will publish an explanatory note later this year providing guidance to users on the quality of current APS and will be used to inform further work ONS is undertaking to improve quality of the survey.`n</p>"
$ws.Range("B10").Value = "July 2025 data"
$ws.Range("C10").Value = "These statistics should be treated as official statistics in development (previously known as experimental statistics). ONS have temporarily paused publication of the snapshot metric for online job adverts so the latest data available covers up to July 2025. More up to date data on new adverts is still available directly from the  <a href=`"https://www.ons.gov.uk/employmentandlabourmarket/peopleinwork/employmentandemployeetypes/datasets/labourdemandvolumesbystandardoccupationclassificationsoc2020uk`">ONS website</a>."
$ws.Range("F10").Value = "<ol>`n  <li>These statistics should be treated as official statistics in development (previously known as experimental statistics), as they are still subject to testing the ability to meet user needs and may be modified in the future.</li>`n<li>Where the same job is identified as being advertised through multiple adverts it is only counted once.</li>`n<li>The method for allocating jobs to occupations (SOC 2020) is based on the job title of the advert and will be developed further in future releases.</li>`n<li>Use caution when interpreting this data. A difference between subgroups does not necessarily imply any causality. There could be other contributing factors at work.</li>`n</ol>"
$ws.Range("I10").Value = "The number of online job adverts in"
$ws.Range("B11").Value = "Mar 2025 data"
$ws.Range("B12").Value = "Dec 2023 - Dec 2024 data"
$ws.Range("B13").Value = "Dec 2023 - Dec 2024 data"
$ws.Range("B14").Value = "AY24/25 data"
$ws.Range("C14").Value = "This indicator shows 19+ further education and skills learner achievements. Further education and skills include apprenticeships and publicly-funded adult learning, including tailored learning, delivered by an FE institution, a training provider or within a local community. "
$ws.Range("D14").Value = "<a href='https://explore-education-statistics.service.gov.uk/data-catalogue/data-set/b930498d-b4f0-416d-a086-7acee1be8179'>Individualised Learner Record</a>"
$ws.Range("E14").Value = "Further education and skills include all age apprenticeships and publicly-funded adult (19+) learning, including tailored learning, delivered by an FE institution, a training provider or within a local community.`nFE and skills does not includer higher education, unless delivered as part of an apprenticeship programme.`nApprenticeships are paid jobs that incorporate on-the-job and off-the-job training leading to nationally recognised qualifications.`nTailored learning is primarily non-qualification based provision that is tailored to the skills needs of the learners, employers and local communities.`nAchievements are the number of learners who successfully complete an individual aim in an academic year."
$ws.Range("I14").Value = "The number of FE achievements in"
$ws.Range("B15").Value = "AY24/25 data"
$ws.Range("C15").Value = "This indicator shows 19+ further education and skills learner achievements. Further education and skills include apprenticeships and publicly-funded adult learning, including tailored learning, delivered by an FE institution, a training provider or within a local community. "
$ws.Range("D15").Value = "<a href='https://explore-education-statistics.service.gov.uk/data-catalogue/data-set/b930498d-b4f0-416d-a086-7acee1be8179'>Individualised Learner Record</a>"
$ws.Range("E15").Value = "Further education and skills include all age apprenticeships and publicly-funded adult (19+) learning, including tailored learning, delivered by an FE institution, a training provider or within a local community.`nFE and skills does not includer higher education, unless delivered as part of an apprenticeship programme.`nApprenticeships are paid jobs that incorporate on-the-job and off-the-job training leading to nationally recognised qualifications.`nTailored learning is primarily non-qualification based provision that is tailored to the skills needs of the learners, employers and local communities.`nAchievements are the number of learners who successfully complete an individual aim in an academic year."
$ws.Range("I15").Value = "FE participation in"
$ws.Range("B16").Value = "AY24/25 data"
$ws.Range("C16").Value = "This indicator shows 19+ further education and skills achievements rate per 100,000 population. Further education and skills include apprenticeships and publicly-funded adult learning, including tailored learning, delivered by an FE institution, a training provider or within a local community. `n"
$ws.Range("D16").Value = "<a href='https://explore-education-statistics.service.gov.uk/data-catalogue/data-set/b930498d-b4f0-416d-a086-7acee1be8179'>Individualised Learner Record</a>"
$ws.Range("E16").Value = "FE and skills does not includer higher education, unless delivered as part of an apprenticeship programme.`nApprenticeships are paid jobs that incorporate on-the-job and off-the-job training leading to nationally recognised qualifications.`nTailored learning is primarily non-qualification based provision that is tailored to the skills needs of the learners, employers and local communities.`nAchievements are the number of learners who successfully complete an individual aim in an academic year."
$ws.Range("F16").Value = "<ol>`n  <li>Achievements included are learners that completed their qualification at any point during the stated academic period.</li>`n <li>Learners achieving more than one course will appear only once in totals.</li>`n <li>Years shown represent academic years.</li>`n<li>Use caution when interpreting this data. A difference between subgroups does not necessarily imply any causality. There could be other contributing factors at work.</li>`n</ol>"
$ws.Range("I16").Value = "The FE achievement rate per 100,000 in"
$ws.Range("B17").Value = "AY24/25 data"
$ws.Range("C17").Value = "This indicator shows 19+ further education and skills participation rate per 100,000 population. Further education and skills include apprenticeships and publicly-funded adult learning, including tailored learning, delivered by an FE institution, a training provider or within a local community. `n"
$ws.Range("D17").Value = "<a href='https://explore-education-statistics.service.gov.uk/data-catalogue/data-set/b930498d-b4f0-416d-a086-7acee1be8179'>Individualised Learner Record</a>"
$ws.Range("E17").Value = "FE and skills does not includer higher education, unless delivered as part of an apprenticeship programme.`nApprenticeships are paid jobs that incorporate on-the-job and off-the-job training leading to nationally recognised qualifications.`nTailored learning is primarily non-qualification based provision that is tailored to the skills needs of the learners, employers and local communities.`nAchievements are the number of learners who successfully complete an individual aim in an academic year."
$ws.Range("F17").Value = "<ol>`n  <li>Participation includes learners that participated at any point during the stated academic period.</li>`n <li>Learners participating in more than one course will appear only once in the grand total.</li>`n <li>Years shown represent academic years.</li>`n<li>Use caution when interpreting this data. A difference between subgroups does not necessarily imply any causality. There could be other contributing factors at work.</li>`n</ol>"
$ws.Range("I17").Value = "The FE participation rate per 100,000 in"
$ws.Range("B18").Value = "Jan-Dec 2024 data"
$ws.Range("I18").Value = "The proportion of people qualified at Level 3 or above in"
$ws.Range("B19").Value = "Jan-Dec 2024 data"
$ws.Range("I19").Value = "The proportion of people qualified at Level 4 or above in"
$ws.Range("B20").Value = "AY23/24 data"
$ws.Range("C20").Value = "Destination measures show the percentage of students going to or remaining in an education, apprenticeship or employment destination in the academic year after completing Key Stage 4 studies (usually aged between 14 to 16). The cohort of learners used in the metrics here completed in AY22/23."
$ws.Range("D20").Value = "<a href = 'https://explore-education-statistics.service.gov.uk/find-statistics/key-stage-4-destination-measures/2022-23'>Key stage 4 destination measures</a>"
$ws.Range("E20").Value = "Destination measures show the percentage of students going to or remaining in an education, apprenticeship or employment destination in the academic year after completing Key Stage 4 studies (usually aged between 14 to 16). The cohort of learners used in the metrics here completed in AY22/23.`nA sustained destination is a count of young people recorded as having sustained participation (education and employment) for a 6 month period in the destination year.`nThis means attending for all of the first two terms of the academic year (e.g. October 2023 to March 2024) at one or more education providers; spending 5 of the 6 months in employment or a combination of the two.`nA sustained apprenticeship is recorded when 6 months continuous participation is recorded at any point in the destination year (between August 2023 and July 2024).`nNot recorded includes pupils who were captured in the destination source data but who failed to meet the sustained participation criteria.`nUnknown (activity not captured): The student was not found to have any participation in education, apprenticeship or employment nor recorded as receiving out-of-work benefits at any point in the year. This also includes not being recorded by their Local Authority as NEET (not engaged in education, employment or training)."
$ws.Range("B21").Value = "AY23/24 data"
$ws.Range("C21").Value = "Destination measures show the percentage of students going to or remaining in an education, apprenticeship or employment destination in the academic year after completing Key Stage 5 studies (usually aged 18). The cohort of learners used in the metrics here completed in AY22/23."
$ws.Range("D21").Value = "<a href = 'https://explore-education-statistics.service.gov.uk/find-statistics/16-18-destination-measures'>16-18 destination measures</a>"
$ws.Range("E21").Value = "Destination measures show the percentage of students going to or remaining in an education, apprenticeship or employment destination in the academic year after completing Key Stage 5 studies (usually aged 18). The cohort of learners used in the metrics here completed in AY22/23.`nA sustained destination is a count of young people recorded as having sustained participation (education and employment) for a 6 month period in the destination year.`nThis means attending for all of the first two terms of the academic year (e.g. October 2023 to March 2024) at one or more education providers; spending 5 of the 6 months in employment or a combination of the two.`nA sustained apprenticeship is recorded when 6 months continuous participation is recorded at any point in the destination year (between August 2023 and July 2024).`nNot recorded includes pupils who were captured in the destination source data but who failed to meet the sustained participation criteria.`nUnknown (activity not captured): The student was not found to have any participation in education, apprenticeship or employment nor recorded as receiving out-of-work benefits at any point in the year. This also includes not being recorded by their Local Authority as NEET (not engaged in education, employment or training)."
$ws.Range("F21").Value = "<ol>`n  <li>Data based on destinations of state-funded mainstream schools and colleges.</li>`n <li>There is no double counting across destinations, a young person is reported in one destination category only.</li>`n <li>If a student is registered as being in education and an apprenticeship, it is recorded as a sustained education and if a student is registered in employment along with an apprenticeship or in education, it is recorded as sustained employment.</li>`n<li>Use caution when interpreting this data. A difference between subgroups does not necessarily imply any causality. There could be other contributing factors at work.</li>`n</ol>"
$ws.Range("B22").Value = "Growth from 2024 to 2035. The LSIP boundaries have changed and new CAs have been created since this data was published so some areas no longer have data"
$ws.Range("E22").Value = "Data presented is the baseline projection. Alternative scenarios are available within the published data for UK only."
$ws.Range("F22").Value = "<ol>`n  <li>The projections presented in this Workbook are calculated from a number of different data sources, using a variety of econometric and statistical techniques. As a result, precise margins of error cannot be assigned to the estimates. For further details, see the Technical Report. </li>`n <li>Industries are based on SIC 2007 codes. </li>`n <li>Time series of the breakdowns can be downloaded in the data download section or in the publication. </li>`n <li>Further breakdowns are available in the published data eg gender, full-time/part-time, as well combined breakdowns. Replacement demand is also available. </li>`n <li>The projections were largely prepared before the Russian invasion of Ukraine, its subsequent economic effects and policy responses to these events. These factors will have a significant impact on the short-term prospects for the British economy and labour market </li>`n<li>Use caution when interpreting this data. A difference between subgroups does not necessarily imply any causality. There could be other contributing factors at work.</li>`n</ol>"
$ws.Range("J22").Value = "Growth from 2024 to 2035"
$ws.Range("B23").Value = "AY22/23 data"
$ws.Range("B24").Value = "AY22/23 data"

# --- Row height changes for rows 14 and 15 (363 -> 334) ---
$ws.Rows.Item(14).RowHeight = 334
$ws.Rows.Item(15).RowHeight = 334

# --- Update selection / scroll position ---
$win = $excel.ActiveWindow
$win.ScrollRow = 14
$win.ScrollColumn = 1
$ws.Range("E16").Select()
